$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.231987118721008
$ws.Range("B1").Value = 3.817236423492432
$ws.Range("C1").Value = 3.51042628288269
$ws.Range("D1").Value = 3.622545003890991
$ws.Range("E1").Value = 1.115609407424927
